# Fixed query issue for C3DC phs002599
#
# The "Treatment Tab" query (row 5, column B) wrapped its REPLACE() call in a
# redundant CONCAT() — i.e. CONCAT(REPLACE(trt.treatment_agent, ';', ', ')).
# Strip the no-op CONCAT wrapper so the cell reads:
#   REPLACE(trt.treatment_agent, ';', ', ') AS "Treatment Agent"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$treatmentCell = $ws.Range("B5")

$oldQuery = $treatmentCell.Value2
$badFragment  = "CONCAT(REPLACE(trt.treatment_agent, ';', ', '))"
$fixedFragment = "REPLACE(trt.treatment_agent, ';', ', ')"

if ($oldQuery.Contains($badFragment)) {
    $newQuery = $oldQuery.Replace($badFragment, $fixedFragment)
    $treatmentCell.Value = $newQuery
}

# The cell's font was nudged down to 11pt (from 12pt used by the other query
# cells) while the fix was made.
$treatmentCell.Font.Size = 11
$treatmentCell.Font.ThemeColor = 1

# Leave the selection/view on the row that was edited.
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$ws.Range("C5").Select() | Out-Null
